# Fehler in Präsentation korrigiert
# Slide 10: center-align the four file-operation captions, widen/move the
# "Insert gesetzliche Krankenkasse.xlsx" label and fix two typo'd captions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# EMU -> point conversion factor (Shape.Left/Width etc. are expressed in points)
$emuPerPoint = 12700

# 1) "Delete Personalnummer.xlsx" -> center align paragraph
$sh1 = $s.Shapes.Item("Textfeld 17")
$sh1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 2) "Insert Mitarbeiter.xlsx" -> center align paragraph
$sh2 = $s.Shapes.Item("Textfeld 18")
$sh2.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 3) "Insert ges. Krankenkasse.xlsx" -> reposition/resize, center align,
#    and fix the abbreviated word in the caption
$sh3 = $s.Shapes.Item("Textfeld 19")
$sh3.Left = 500656 / $emuPerPoint
$sh3.Width = 1983112 / $emuPerPoint
$sh3.TextFrame.TextRange.Text = "Insert gesetzliche Krankenkasse.xlsx"
$sh3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 4) "Update Krankenversicherungsbeitraege.xlsx" -> center align and correct caption
$sh4 = $s.Shapes.Item("Textfeld 20")
$sh4.TextFrame.TextRange.Text = "Update adresse.xlsx"
$sh4.TextFrame.TextRange.ParagraphFormat.Alignment = 2
